$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1-3: update the first three rows' text
$t.Cell(1,1).Range.Text = "0M"
$t.Cell(2,1).Range.Text = "0M"
$t.Cell(3,1).Range.Text = "0M"

# 4: insert 10 new rows before the current row 4 (still holding "0"),
#    carrying the per-iteration timing values that used to be packed
#    (tab-separated) into the three rows near the end of the table.
# Each Rows.Add(beforeRow) inserts immediately above beforeRow, so walk the
# values in reverse to end up with them in forward reading order.
$newValues = @("21", "0.00002", "0.00005", "0.00003", "0.00000", "0.00003", "0.00004", "0.00004", "0.00068", "100.0")
$beforeRow = $t.Rows.Item(4)
for ($i = $newValues.Length - 1; $i -ge 0; $i--) {
    $newRow = $t.Rows.Add($beforeRow)
    $newRow.Cells.Item(1).Range.Text = $newValues[$i]
    $beforeRow = $newRow
}

# 5-7: the three rows that used to hold tab-separated per-iteration data
#      (now shifted down by the 10 inserted rows) collapse back down to a
#      single summary value each.
$t.Cell(44,1).Range.Text = "100"
$t.Cell(45,1).Range.Text = "0"
$t.Cell(46,1).Range.Text = "27"
